# B6-PowerPoint.pptx edit: (1) switch the three tables on slides 14-16 to the
# new table style, (2) swap the two theme colour palettes ("Integral" /
# "Office Theme") that the deck's slide master and notes master use.

$p = $ppt.ActivePresentation

# --- 1) Table style id change -------------------------------------------------
# Three tables (slide 14, 15, 16 - each the first shape on its slide) move
# from the default "Table_0" style to the built-in style
# {8170544A-4285-4BE2-A49C-0E0B1320F64B}.
$newTableStyle = "{8170544A-4285-4BE2-A49C-0E0B1320F64B}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle, $false)
        }
    }
}

# --- 2) Theme colour swap -----------------------------------------------------
# The deck carries two colour themes: the slide master uses "Integral" (Red
# Violet accents) while the notes master uses "Office Theme". The edit swaps
# which palette lives in which theme part, so the slide master should end up
# with the classic "Office Theme" colours.
$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
